$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "from manju branch"
$ws.Range("B4").Select()
